$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 102.833336
$ws.Range("I33").Value = 93.40000000000001
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 93.40000000000001
$ws.Range("L33").Value = 150
$ws.Range("M33").Value = 135.6
$ws.Range("N33").Value = -608
# Row 40
$ws.Range("H40").Value = 2000
$ws.Range("J40").Value = 2000
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2350
# Row 51
$ws.Range("H51").Value = 8382.5
$ws.Range("I51").Value = 8449
$ws.Range("J51").Value = 8249.5
$ws.Range("K51").Value = 8449
$ws.Range("L51").Value = 8249.5
$ws.Range("M51").Value = -7965
$ws.Range("N51").Value = -9217.5
# Row 53
$ws.Range("H53").Value = 826.25
$ws.Range("I53").Value = 601.6667
$ws.Range("K53").Value = 601.6667
$ws.Range("M53").Value = 35.33330000000001
# Row 58
$ws.Range("H58").Value = 1128
$ws.Range("J58").Value = 1499
$ws.Range("L58").Value = 4497
$ws.Range("N58").Value = -4797
# Row 64
$ws.Range("H64").Value = 7000
$ws.Range("I64").Value = 4000
$ws.Range("K64").Value = 4000
$ws.Range("M64").Value = -3752
# Row 67
$ws.Range("H67").Value = 7000
$ws.Range("I67").Value = 4000
$ws.Range("K67").Value = 4000
$ws.Range("M67").Value = -3142
# Row 82
$ws.Range("H82").Value = 735
$ws.Range("I82").Value = 735
$ws.Range("K82").Value = 2205
$ws.Range("M82").Value = -1799
# Row 85
$ws.Range("H85").Value = 735
$ws.Range("I85").Value = 735
$ws.Range("K85").Value = 2205
$ws.Range("M85").Value = -801
# Row 127
$ws.Range("H127").Value = 5339.4
$ws.Range("I127").Value = 6049.25
$ws.Range("J127").Value = 2500
$ws.Range("K127").Value = 18147.75
$ws.Range("L127").Value = 7500
$ws.Range("M127").Value = -13187.75
$ws.Range("N127").Value = -17420
# Row 129
$ws.Range("H129").Value = 14291765
$ws.Range("I129").Value = 50001100
$ws.Range("J129").Value = 8032
$ws.Range("K129").Value = 150003300
$ws.Range("L129").Value = 24096
$ws.Range("M129").Value = -149998300
$ws.Range("N129").Value = -34096
# Row 132
$ws.Range("H132").Value = 2124.5
$ws.Range("I132").Value = 2124.5
$ws.Range("K132").Value = 6373.5
$ws.Range("M132").Value = -3843.5

$ws = $wb.Worksheets.Item("ARM")
# Row 55
$ws.Range("H55").Value = 60000
$ws.Range("J55").Value = 60000
$ws.Range("L55").Value = 60000
$ws.Range("N55").Value = -60630
# Row 115
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -43134

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 8999.666999999999
$ws.Range("I86").Value = 3499
$ws.Range("K86").Value = 3499
$ws.Range("M86").Value = -2376
# Row 88
$ws.Range("H88").Value = 9933.333000000001
$ws.Range("J88").Value = 9933.333000000001
$ws.Range("L88").Value = 9933.333000000001
$ws.Range("N88").Value = -10745.333
# Row 89
$ws.Range("H89").Value = 8999.666999999999
$ws.Range("I89").Value = 3499
$ws.Range("K89").Value = 17495
$ws.Range("M89").Value = -11879
# Row 91
$ws.Range("H91").Value = 9933.333000000001
$ws.Range("J91").Value = 9933.333000000001
$ws.Range("L91").Value = 9933.333000000001
$ws.Range("N91").Value = -12741.333
# Row 95
$ws.Range("H95").Value = 15816.333
$ws.Range("J95").Value = 15816.333
$ws.Range("L95").Value = 15816.333
$ws.Range("N95").Value = -21308.333
# Row 105
$ws.Range("H105").Value = 1932.6666
$ws.Range("I105").Value = 1932.6666
$ws.Range("K105").Value = 1932.6666
$ws.Range("M105").Value = -185.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 56.6
$ws.Range("I7").Value = 58.25
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 58.25
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = 54.75
$ws.Range("N7").Value = -276
# Row 22
$ws.Range("H22").Value = 241.57143
$ws.Range("I22").Value = 248.5
$ws.Range("J22").Value = 232.33333
$ws.Range("K22").Value = 248.5
$ws.Range("L22").Value = 232.33333
$ws.Range("M22").Value = 101.5
$ws.Range("N22").Value = -932.3333299999999
# Row 57
$ws.Range("H57").Value = 48999.5
$ws.Range("J57").Value = 48999.5
$ws.Range("L57").Value = 48999.5
$ws.Range("N57").Value = -50119.5
# Row 86
$ws.Range("H86").Value = 10938.75
$ws.Range("I86").Value = 10700.6
$ws.Range("K86").Value = 10700.6
$ws.Range("M86").Value = -9577.6
# Row 89
$ws.Range("H89").Value = 10938.75
$ws.Range("I89").Value = 10700.6
$ws.Range("K89").Value = 53503
$ws.Range("M89").Value = -47887
# Row 134
$ws.Range("H134").Value = 7000
$ws.Range("I134").Value = 7000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 21000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -18465
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 1
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 109
# Row 34
$ws.Range("H34").Value = 538
$ws.Range("I34").Value = 575
$ws.Range("J34").Value = 501
$ws.Range("K34").Value = 1725
$ws.Range("L34").Value = 1503
$ws.Range("M34").Value = -1641
$ws.Range("N34").Value = -1671
# Row 81
$ws.Range("H81").Value = 500
$ws.Range("J81").Value = 500
$ws.Range("L81").Value = 1500
$ws.Range("N81").Value = -3746
# Row 84
$ws.Range("H84").Value = 500
$ws.Range("J84").Value = 500
$ws.Range("L84").Value = 4500
$ws.Range("N84").Value = -15732
# Row 141
$ws.Range("H141").Value = 2769
$ws.Range("I141").Value = 2769
$ws.Range("K141").Value = 8307
$ws.Range("M141").Value = -3127

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5875
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730
# Row 73
$ws.Range("H73").Value = 5875
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064
# Row 80
$ws.Range("H80").Value = 3319.4
$ws.Range("J80").Value = 1750
$ws.Range("L80").Value = 1750
$ws.Range("N80").Value = -3746
# Row 83
$ws.Range("H83").Value = 3319.4
$ws.Range("J83").Value = 1750
$ws.Range("L83").Value = 8750
$ws.Range("N83").Value = -18734

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 1065.6666
$ws.Range("I100").Value = 99.5
$ws.Range("J100").Value = 2998
$ws.Range("K100").Value = 99.5
$ws.Range("L100").Value = 2998
$ws.Range("M100").Value = 441.5
$ws.Range("N100").Value = -4080
# Row 122
$ws.Range("H122").Value = 4992.3335
$ws.Range("I122").Value = 5551
$ws.Range("J122").Value = 3875
$ws.Range("K122").Value = 16653
$ws.Range("L122").Value = 11625
$ws.Range("M122").Value = -14203
$ws.Range("N122").Value = -16525

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 986.6
$ws.Range("I81").Value = 986.6
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1973.2
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -912.2
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 986.6
$ws.Range("I84").Value = 986.6
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9866
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4562
$ws.Range("N84").ClearContents()
# Row 122
$ws.Range("H122").Value = 4500
$ws.Range("I122").Value = 4500
$ws.Range("K122").Value = 13500
$ws.Range("M122").Value = -11050
# Row 126
$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10000
$ws.Range("K126").Value = 30000
$ws.Range("M126").Value = -27530
# Row 132
$ws.Range("H132").Value = 9000
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060

